# Update Work Week and Social Spending
# Revises the Pakistan GDP-per-Capita series (rows for years 1950-2010) to the
# updated Clio-Infra figures, and appends six new years (2011-2016).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Updated "Data" (GDP per Capita) values for years 1950..2016, row 2..68.
# (Years 1950-2010 were revised; 2011-2016 are brand new rows.)
$values = @(
    "1025",
    "969",
    "950",
    "1015",
    "1014",
    "1012",
    "1017",
    "1036",
    "1025",
    "1009",
    "1031",
    "1066",
    "1114",
    "1152",
    "1208",
    "1229",
    "1294",
    "1307",
    "1361",
    "1411",
    "1517",
    "1484",
    "1455",
    "1521",
    "1533",
    "1559",
    "1604",
    "1631",
    "1720",
    "1733",
    "1851",
    "1924",
    "2002",
    "2077",
    "2106",
    "2230",
    "2302",
    "2364",
    "2442",
    "2487",
    "2531",
    "2644.57421556903",
    "2732.38162109905",
    "2777.33203198315",
    "2855.90619647847",
    "2947.68262621299",
    "2979.80822550806",
    "2986.70556564547",
    "3035.45170157769",
    "3111.95857036277",
    "3154.51158006955",
    "3176.94677658726",
    "3265.68661537707",
    "3437.9918213429",
    "3664.53343212702",
    "3883.93451491353",
    "4062.66645699341",
    "4160.36516936913",
    "4224.40480371735",
    "4288.46181926384",
    "4354.26695082655",
    "4464",
    "4569",
    "4703",
    "4850",
    "5030",
    "5223"
)

$firstRow = 2
$lastRow = $firstRow + $values.Length - 1

# Column E ("Data") holds these figures as text in the workbook (not numbers),
# so force a Text number format before writing them, otherwise the
# numeric-looking strings would be auto-converted to real numbers.
$ws.Range("E$firstRow`:E$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $values[$i]
}

# Fill in the rest of the six brand-new rows (2011-2016): Country Code,
# Country Name, Indicator and Year columns follow the same pattern as all
# the other rows in the sheet.
$newYears = @(2011, 2012, 2013, 2014, 2015, 2016)
for ($i = 0; $i -lt $newYears.Length; $i++) {
    $row = 63 + $i
    $ws.Cells.Item($row, 1).Value = 586
    $ws.Cells.Item($row, 2).Value = "Pakistan"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $newYears[$i]
}
